{"js": "// The abstract paragraph is being reworded: the sentences describing where\n// the number-plate-detection system is used, and how it helps, are rewritten\n// with clearer/expanded phrasing (e.g. fastag example, \"keep a track of\",\n// \"primary objective\", \"resulting images\", \"segmentation processes\",\n// \"Optical character reading\", \"simulated and implemented\").\n//\n// The unchanged lead-in sentence (\"The NPD (Number Plate Detection) system is\n// an image processing technology.\") stays untouched, as does the\n// \"system inputs a vehicle image and separates the number plate from the\n// rest of the image\" clause in the middle of the paragraph. Everything else\n// in the rest of the paragraph is reworded.\n\nconst body = context.document.body;\n\n// --- First block: from \"This type of system...\" through to \"...The\" (right\n// before the existing \" system inputs a vehicle image...\" text). ---\nconst oldBlock1 =\n  \" This type of system is widely used in Traffic control areas, tolling, parking area etc. \" +\n  \" It is one of the necessary systems designed to detect the vehicle number plate. With the development of this system, it becomes easy to keep a record\" +\n  \" of the increasing number of vehicles and\" +\n  \" use it whenever required. The main objective is to design an efficient vehicle identification system by using vehicle number plate. The\";\n\nconst newBlock1 =\n  \" \" +\n  \"Such systems\" +\n  \" \" +\n  \"are\" +\n  \" widely used in Traffic control areas, tolling\" +\n  \" \" +\n  \"(for e\" +\n  \".\" +\n  \"g.\" +\n  \",\" +\n  \" fastag)\" +\n  \", parking area etc.\" +\n  \" It\" +\n  \" becomes easy to \" +\n  \"keep a track of \" +\n  \"the increasing number of vehicles and\" +\n  \" use it whenever required\" +\n  \" with the development of this system\" +\n  \". The \" +\n  \"primary\" +\n  \" objective is to design \" +\n  \"a \" +\n  \"vehicle identification system by using \" +\n  \"the\" +\n  \" number plate\" +\n  \" of the vehicle\" +\n  \". The\";\n\n// --- Second block: from \". The captured images...\" through to the end of\n// the paragraph. ---\nconst oldBlock2 =\n  \". The captured images are then extracted by using the segmentation process. Optical character recognition is used to identify the characters. The system is implemented and simulated on MATLAB and performance is tested on real images. \";\n\nconst newBlock2 =\n  \". The \" +\n  \"resulting\" +\n  \" images are\" +\n  \" \" +\n  \"extracted by using the segmentation process\" +\n  \"es\" +\n  \". Optical character re\" +\n  \"ading is used to convert images to machine encoded text\" +\n  \". The system is \" +\n  \"simulated\" +\n  \" and \" +\n  \"implemented\" +\n  \" on MATLAB and performance is tested on real images. \";\n\nasync function replaceOnce(searchText, replacementText) {\n  const results = body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not locate expected text: \" + searchText);\n  }\n\n  results.items[0].insertText(replacementText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Replace the tail block first so the offsets/text used to find block 1\n// remain valid (the two blocks don't overlap, but doing the later one first\n// keeps this resilient either way).\nawait replaceOnce(oldBlock2, newBlock2);\nawait replaceOnce(oldBlock1, newBlock1);\n", "ps1": "# The abstract paragraph is being reworded: the sentences describing where\n# the number-plate-detection system is used, and how it helps, are rewritten\n# with clearer/expanded phrasing (e.g. fastag example, \"keep a track of\",\n# \"primary objective\", \"resulting images\", \"segmentation processes\",\n# \"Optical character reading\", \"simulated and implemented\").\n#\n# The unchanged lead-in sentence (\"The NPD (Number Plate Detection) system is\n# an image processing technology.\") stays untouched, as does the\n# \"system inputs a vehicle image and separates the number plate from the\n# rest of the image\" clause in the middle of the paragraph. Everything else\n# in the rest of the paragraph is reworded.\n\n$d = $word.ActiveDocument\n\n# --- First block: from \"This type of system...\" through to \"...The\" (right\n# before the existing \" system inputs a vehicle image...\" text). ---\n$oldBlock1 = \" This type of system is widely used in Traffic control areas, tolling, parking area etc. \" + `\n  \" It is one of the necessary systems designed to detect the vehicle number plate. With the development of this system, it becomes easy to keep a record\" + `\n  \" of the increasing number of vehicles and\" + `\n  \" use it whenever required. The main objective is to design an efficient vehicle identification system by using vehicle number plate. The\"\n\n$newBlock1 = \" Such systems are widely used in Traffic control areas, tolling (for e.g., fastag), parking area etc.\" + `\n  \" It becomes easy to keep a track of the increasing number of vehicles and use it whenever required with the development of this system.\" + `\n  \" The primary objective is to design a vehicle identification system by using the number plate of the vehicle. The\"\n\n# --- Second block: from \". The captured images...\" through to the end of\n# the paragraph. ---\n$oldBlock2 = \". The captured images are then extracted by using the segmentation process. Optical character recognition is used to identify the characters. The system is implemented and simulated on MATLAB and performance is tested on real images. \"\n\n$newBlock2 = \". The resulting images are extracted by using the segmentation processes. Optical character reading is used to convert images to machine encoded text. \" + `\n  \"The system is simulated and implemented on MATLAB and performance is tested on real images. \"\n\n# wdReplaceOne = 1\n$wdReplaceOne = 1\n\n$rng1 = $d.Content\n$found1 = $rng1.Find.Execute($oldBlock1, $false, $false, $false, $false, $false, $true, 1, $false, $newBlock1, $wdReplaceOne)\nif (-not $found1) {\n  throw \"Could not locate the first target block of abstract text to replace.\"\n}\n\n$rng2 = $d.Content\n$found2 = $rng2.Find.Execute($oldBlock2, $false, $false, $false, $false, $false, $true, 1, $false, $newBlock2, $wdReplaceOne)\nif (-not $found2) {\n  throw \"Could not locate the second target block of abstract text to replace.\"\n}\n"}
